$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.949.17'
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").Value = '1.777.72'
$ws.Range("E3").Value = '  +0.41%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.03%  '

$ws.Range("E6").Value = '  +1.38%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.12'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.290'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0698'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0938'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.99%  '

$ws.Range("D12").Value = '2.036.00'
$ws.Range("E12").Value = '  +0.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.92%  '

$ws.Range("D14").Value = '1.771.87'
$ws.Range("E14").Value = '  +0.11%  '

$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '33.956.90'
$ws.Range("E15").Value = '  +0.06%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.619'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.73%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.16%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.33%  '

$ws.Range("D20").Value = '0.0₃0780'
$ws.Range("E20").Value = '  +0.79%  '

$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.67'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.97%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.52'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.63%  '

$ws.Range("E28").Value = '  +1.19%  '

$ws.Range("E29").Value = '  +0.21%  '

$ws.Range("E31").Value = '  -1.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.64'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.65%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.60%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.47%  '

$ws.Range("D35").Value = '1.390.82'
$ws.Range("E35").Value = '  -0.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.658'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.27%  '

$ws.Range("E37").Value = '  -0.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0186'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.51%  '

$ws.Range("E39").Value = '  +6.59%  '

$ws.Range("E40").Value = '  +0.88%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.911'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.74%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '77.88'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.79%  '

$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").Value = '0.0₆0145'
$ws.Range("E44").Value = '  +21.49%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.34'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +12.85%  '

$ws.Range("E46").Value = '  +4.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '108.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0497'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.78%  '

$ws.Range("E49").Value = '  -0.41%  '

$ws.Range("D50").Value = '1.935.67'
$ws.Range("E50").Value = '  +0.92%  '

$ws.Range("E51").Value = '  +0.50%  '
